$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force Text number format before writing a value so that
# numeric-looking strings (e.g. "26.656.91", "0.0620", "15.50") are stored
# verbatim as text instead of being auto-converted/rounded by Excel, then
# clear the formatting again so the cell keeps its original (default) style
# -- only the text content changes, matching the source diff exactly.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.656.91"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.632.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0620"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.96"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.858.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.658.93"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.93"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.634.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Dai"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.01"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C20").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "211.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Chainlink"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C22").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Toncoin"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C23").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -8.31%  "
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.03%  "
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.256.76"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "MXToken"
$ws.Range("B42").ClearFormats()
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C42").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("B43").ClearFormats()
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C43").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.769.77"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.12"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.65"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.406"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0954"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.66%  "
$ws.Range("E51").ClearFormats()
